$d = $word.ActiveDocument

$replacements = @(
    @("Deskribí kiko bo yu ta hasiendo.", "Deskribí kiko bo yu ta hasi."),
    @("“Bo ta trahando un toren ku blòki hel i ata akí e kòrá ku bo ta poniendo riba dje. Awor bo ta hinka e sinku bestianan den e garoshi i nan ta bai di biahe.” ", "“Bo ta traha un toren ku blòki hel i aki tin e kòrá ku bo ta pone riba dje. Awor bo ta hinka e sinku bestianan den e garoshi i nan ta bai keiru.” "),
    @("“Mi ta mira ku bo ta trahando duru riba bo hùiswèrk. Mi ta mira ku bo ta purbando duru pa solushoná e problema akí. Sigui asina!” ", "“Mi ta mira ku bo ta traha duru di bèrda riba bo hùiswèrk. Mi ta mira ku bo ta hasi bo bèst pa solushoná e problema akí. Sigui asina!” "),
    @("Esaki por sinti straño na promé instante ya ku hopi mayor/edukadónan no ta kustumá di papia ku nan yunan durante ora di wega. Sinembargo, ku práktika, e ta bira mas fásil – meskos ku tur otro kos! ", "Esaki por sinti straño na promé instante ya ku hopi mayor/edukadónan no ta kustumá di papia ku nan yunan ora ku nan ta hunga. Sinembargo, ku práktika, e ta bira mas fásil – meskos ku tur otro kos! "),
    @("Ounke mayornan/edukadónan hopi bia ta bisa ku deskribí kiko nan yu ta hasiendo ta pone nan sinti nan mes inkómodo, pa nan sorpresa, nan yunan por lo general no ta nota ku nan ta inkómodo i e yunan ta respondé hopi positivo na e atenshon akí.  ", "Ounke mayornan/edukadónan hopi bia ta bisa ku deskribí kiko nan yu ta hasi ta pone nan sinti nan mes inkómodo, pa nan sorpresa, nan yunan por lo general no ta nota ku nan ta inkómodo i e yunan ta respondé hopi positivo riba e atenshon akí.  "),
    @("Ta normal pa mayornan/edukadónan kuminsá ku hasi pregunta en bes di “Bisa Loke Bo Ta Mira”. Ta bo ròl komo fasilitadó pa yuda nan siña kon pa simplemente deskribí kiko e mucha ta hasiendo en bes di hasi pregunta. ", "Ta normal pa mayornan/edukadónan kuminsá ku hasi pregunta na lugá di “Bisa Loke Bo Ta Mira”. Ta bo ròl komo fasilitadó pa yuda nan siña kon pa simplemente deskribí kiko e mucha ta hasi na lugá di hasi pregunta. "),
    @("Kòrda, esaki ta e wega di e mucha. No tin bon ni malu den wega i denter di loke ta rasonabel. Loke un mucha skohe pa hunga ta bon p'e. E mayor su trabou ta pa mustra interes i bisa algu bunita. ", "Kòrda, esaki ta e mucha su wega. No tin bon ni malu den wega i denter di loke ta rasonabel. Loke un mucha skohe pa hunga ta bon p'e. E mayor su trabou ta pa mustra interes i bisa algu bunita. "),
    @("Durante Tempu huntu ku nan adolesentenan, mayornan/edukadónan tambe por permití nan yunan tuma e liderazgo. E hóbennan por skohe tokante kiko e aktividat òf kombersashon ta bai. E por asta nifiká ku e mayornan/edukadónan i hóbennan por pasa Tempu Huntu sintá banda di otro ta hasi aktividatnan paralelo hasiendo opservashonnan de bes en kuando. Ta importante pa enkurashá mayornan/edukadónan pa permití nan yunan papia tokante e kosnan ku ta importante pa nan.  ", "Durante Tempu pa Abo ku bo Yu adolesentenan, mayornan/edukadónan tambe por permití nan yunan tuma e liderazgo. E hóbennan por skohe e aktividat òf kombersashon. E por asta nifiká ku e mayornan/edukadónan i hóbennan por pasa Tempu Huntu sintá banda di otro ta hasi aktividatnan paralelo hasiendo opservashonnan de bes en kuando. Ta importante pa enkurashá mayornan/edukadónan pa permití nan yunan papia tokante e kosnan ku ta importante pa nan.  "),
    @("Un otro bon manera ku mayornan/edukadónan por desaroyá relashonnan positivo durante Tempu ku nan yunan ta pa pasa tempu ku nan ora nan ta hasiendo algu ku nan ta gusta. Por ehèmpel, mira nan yunan hunga un partido deportivo òf mustra un moveshon nobo di baile. ", "Un otro bon manera ku mayornan/edukadónan por desaroyá relashonnan positivo durante Tempu pa Abo ku bo Yu ta pa pasa tempu ku nan ora nan ta hasi algu ku nan ta gusta. Por ehèmpel, mira nan yunan hunga un partido deportivo òf mustra un moveshon nobo di baile. "),
    @("Ora bo ta papia ku un mucha ku desabilidat, kòrda ku bo ta interkambiando ku un mucha, Tene e mesun tono i idioma ku bo lo hasi ku kualke mucha di un edat similar. Si e mucha no por tende, ta importante pa bo sinta ketu, inkluso bo kabes, miéntras bo ta papia ku e mucha. Sòru pa bo wak nan ora bo ta papia i sòru pa nan wak bo i por mira bo kara i boka.  ", "Ora bo ta papia ku un mucha ku desabilidat, kòrda ku bo ta interkambiando ku un mucha, Tene e mésun tono i idioma ku bo lo hasi ku kualke mucha di un edat similar. Si e mucha no por tende, ta importante pa bo sinta ketu, inkluso bo kabes, miéntras bo ta papia ku e mucha. Sòru pa bo wak nan ora bo ta papia i sòru pa nan wak bo i por mira bo kara i boka.  "),
    @("Muchanan ku desabilidat ku no ta kustumbrá ku ta puntra nan nan opinion òf ku no ta kustumbrá ku hende ta skucha nan, mester di mas tempu pa krea konfiansa i seguridat. Lo por tuma tempu pa eksplorá e mihó maneranan di komuniká ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di Pasa Tempu huntu ku nan yunan. ", "Muchanan ku desabilidat ku no ta kustumbrá ku ta puntra nan, nan opinion òf ku no ta kustumbrá ku hende ta skucha nan, mester di mas tempu pa krea konfiansa i seguridat. Lo por tuma tempu pa eksplorá e mihó maneranan di komuniká ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di Pasa Tempu Abo ku bo YU ku nan yunan. "),
    @("Permití muchanan ku desabilidat tuma e mesun riesgonan ku otro muchanan pa asina yuda nan haña konfiansa. Protekshon di mas ta stroba e muchanan di eksplorá, deskubrí kiko ta posibel i siña kon pa mantené nan mes sigur. ", "Permití muchanan ku desabilidat tuma e mésun riesgonan ku otro muchanan pa asina yuda nan haña konfiansa. Protekshon di mas ta stroba e muchanan di eksplorá, deskubrí kiko ta posibel i siña kon pa mantené nan mes sigur. "),
    @("Enfoká riba reforsá e puntonan fuerte i abilidatnan di kada mucha en bes di e kosnan ku nan no por hasi, por ehèmpel, un persona ku ta usa ròlstul por tin brasa i mannan fuerte, un mucha surdu por ta bon den pintamentu", "Enfoká riba reforsá e puntonan fuerte i abilidatnan di kada mucha en bes di e kosnan ku nan no por hasi, por ehèmpel, un persona ku ta usa ròlstul por tin brasa i mannan fuerte, un mucha surdu por ta bon den pintamentu "),
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    $rng = $d.Content.Duplicate
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $rng.Text = $replace
    } else {
        Write-Host "NOT FOUND: $find"
    }
}

Write-Host "Done"
